# Applies the "Trade #36 closed" update to live_trading_results.xlsx
# Sheets (in workbook order): 1=Summary, 2=Strategy Status, 3=All Trades, 4=MarketMaking

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: Summary
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)
$summary.Range("B3").Value = 1400.49   # Current Capital
$summary.Range("B4").Value = 0.29      # Total P&L $
$summary.Range("B5").Value = 0.09      # Total P&L %
$summary.Range("B6").Value = 64        # Total Trades
$summary.Range("B7").Value = 30        # Winning Trades
$summary.Range("B9").Value = 46.88     # Win Rate %

# ---------------------------------------------------------------------------
# Sheet 2: Strategy Status (row 5 = MarketMaking)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item(2)
$status.Range("C5").Value = 100.49     # Capital
$status.Range("D5").Value = 31         # Trades
$status.Range("E5").Value = 0.18       # P&L $
$status.Range("F5").Value = 0.49       # P&L %
$status.Range("G5").Value = 51.61      # Win Rate %

# ---------------------------------------------------------------------------
# Sheet 3: All Trades
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item(3)

# Row 65 (Trade #64) - closes the trade early
$allTrades.Cells.Item(65, 7).Value = 0.86           # G: Exit Price
$allTrades.Cells.Item(65, 8).Value = "CLOSED"       # H: Status
$allTrades.Cells.Item(65, 9).Value = 10.2564        # I: P&L %
$allTrades.Cells.Item(65, 10).Value = 0.08          # J: P&L $
$allTrades.Cells.Item(65, 11).Value = 100.49        # K: Capital After
$allTrades.Cells.Item(65, 12).Value = "early_exit"  # L: Exit Reason
$allTrades.Cells.Item(65, 13).Value = 0.11          # M: Duration (min)

# New row 98 (Trade #97) - newly opened trade appended at the end
$r = 98
$allTrades.Cells.Item($r, 1).Value = 97
$allTrades.Cells.Item($r, 2).NumberFormat = "@"
$allTrades.Cells.Item($r, 2).Value = "2026-02-17"   # B: Date (kept as text)
$allTrades.Cells.Item($r, 3).Value = "21:02:57"     # C: Time
$allTrades.Cells.Item($r, 4).Value = "MarketMaking" # D: Strategy
$allTrades.Cells.Item($r, 5).Value = "DOWN"         # E: Side
$allTrades.Cells.Item($r, 6).Value = 0.78           # F: Entry Price
# G: Exit Price left blank (trade still open) - copy an already-blank cell so it
# still materializes as an explicit (empty) cell instead of being omitted
$allTrades.Range("G66").Copy($allTrades.Range("G98"))
$allTrades.Cells.Item($r, 8).Value = "OPEN"         # H: Status
$allTrades.Cells.Item($r, 9).Value = 0              # I: P&L %
$allTrades.Cells.Item($r, 10).Value = 0             # J: P&L $
$allTrades.Cells.Item($r, 11).Value = 100.4110412885904  # K: Capital After
# L: Exit Reason left blank (trade still open)
$allTrades.Range("L66").Copy($allTrades.Range("L98"))
$allTrades.Cells.Item($r, 13).Value = 0             # M: Duration (min)
$allTrades.Cells.Item($r, 14).Value = 0             # N: Entry Slippage (bps)
$allTrades.Cells.Item($r, 15).Value = 0             # O: Exit Slippage (bps)
$allTrades.Cells.Item($r, 16).Value = 0.6           # P: Confidence
$allTrades.Cells.Item($r, 17).Value = "Normal spread capture: 19600 bps"  # Q: Entry Reason

# ---------------------------------------------------------------------------
# Sheet 4: MarketMaking
# ---------------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item(4)

# Row 32 (Trade #64) - closes the trade early
$marketMaking.Cells.Item(32, 7).Value = 0.86           # G: Exit Price
$marketMaking.Cells.Item(32, 8).Value = "CLOSED"       # H: Status
$marketMaking.Cells.Item(32, 9).Value = 10.2564        # I: P&L %
$marketMaking.Cells.Item(32, 10).Value = 0.08          # J: P&L $
$marketMaking.Cells.Item(32, 11).Value = 100.49        # K: Capital After
$marketMaking.Cells.Item(32, 16).Value = "early_exit"  # P: Exit Reason
$marketMaking.Cells.Item(32, 17).Value = 0.11          # Q: Duration (min)

# New row 65 (Trade #97) - newly opened trade appended at the end
$r2 = 65
$marketMaking.Cells.Item($r2, 1).Value = 97
$marketMaking.Cells.Item($r2, 2).NumberFormat = "@"
$marketMaking.Cells.Item($r2, 2).Value = "2026-02-17"   # B: Date (kept as text)
$marketMaking.Cells.Item($r2, 3).Value = "21:02:57"     # C: Time
$marketMaking.Cells.Item($r2, 4).Value = "MarketMaking" # D: Strategy
$marketMaking.Cells.Item($r2, 5).Value = "DOWN"         # E: Side
$marketMaking.Cells.Item($r2, 6).Value = 0.78           # F: Entry Price
# G: Exit Price left blank (trade still open) - copy an already-blank cell so it
# still materializes as an explicit (empty) cell instead of being omitted
$marketMaking.Range("G33").Copy($marketMaking.Range("G65"))
$marketMaking.Cells.Item($r2, 8).Value = "OPEN"         # H: Status
$marketMaking.Cells.Item($r2, 9).Value = 0              # I: P&L %
$marketMaking.Cells.Item($r2, 10).Value = 0             # J: P&L $
$marketMaking.Cells.Item($r2, 11).Value = 100.4110412885904  # K: Capital After
$marketMaking.Cells.Item($r2, 12).Value = 0             # L: Entry Slippage (bps)
$marketMaking.Cells.Item($r2, 13).Value = 0             # M: Exit Slippage (bps)
$marketMaking.Cells.Item($r2, 14).Value = 0.6           # N: Confidence
$marketMaking.Cells.Item($r2, 15).Value = "Normal spread capture: 19600 bps"  # O: Entry Reason
# P: Exit Reason left blank (trade still open)
$marketMaking.Range("P33").Copy($marketMaking.Range("P65"))
$marketMaking.Cells.Item($r2, 17).Value = 0             # Q: Duration (min)
